$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-02 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("66÷5=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "86÷9=9, 5", 2) | Out-Null
$d.Content.Find.Execute("78÷7=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷2=30, 0", 2) | Out-Null
$d.Content.Find.Execute("24÷8=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "17÷6=2, 5", 2) | Out-Null
$d.Content.Find.Execute("40÷4=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2) | Out-Null
$d.Content.Find.Execute("53÷7=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "76÷8=9, 4", 2) | Out-Null
$d.Content.Find.Execute("99÷6=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=5, 0", 2) | Out-Null
$d.Content.Find.Execute("61÷6=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=13, 1", 2) | Out-Null
$d.Content.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=11, 0", 2) | Out-Null
$d.Content.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("26÷5=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2) | Out-Null
$d.Content.Find.Execute("79÷9=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "18÷5=3, 3", 2) | Out-Null
$d.Content.Find.Execute("16÷3=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=11, 3", 2) | Out-Null
$d.Content.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=4, 5", 2) | Out-Null
$d.Content.Find.Execute("65÷2=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "61÷9=6, 7", 2) | Out-Null
$d.Content.Find.Execute("28÷8=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=10, 1", 2) | Out-Null
$d.Content.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "82÷3=27, 1", 2) | Out-Null
$d.Content.Find.Execute("97÷6=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷5=13, 1", 2) | Out-Null
$d.Content.Find.Execute("37÷6=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=6, 0", 2) | Out-Null
$d.Content.Find.Execute("68÷9=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
$d.Content.Find.Execute("30÷8=3, 6", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=11, 7", 2) | Out-Null
$d.Content.Find.Execute("51÷5=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=15, 0", 2) | Out-Null
$d.Content.Find.Execute("77÷5=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2) | Out-Null
$d.Content.Find.Execute("47÷8=5, 7", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=6, 0", 2) | Out-Null
$d.Content.Find.Execute("18÷7=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2) | Out-Null
$d.Content.Find.Execute("73÷5=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=46, 1", 2) | Out-Null

Write-Host "All replacements applied"
